# finish animations, update sprite inventory sheet
#
# Adds missing W (px) / H (px) values to several existing sprite rows,
# renames a couple of source-file prefixes (grid -> hud/grid,
# hud/arrow -> hud/targeting-arrow), and appends five new "Start Screen"
# sprite rows (Logo, Play Over, Play Out, Play Down, Play Up).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pickups: Box / Diamond now have W/H (px) ---
$ws.Range("F33").Value = 24
$ws.Range("G33").Value = 24

$ws.Range("F34").Value = 24
$ws.Range("G34").Value = 24

# --- FX: Hit / Trail / Swoosh now have W/H (px) ---
$ws.Range("F36").Value = 16
$ws.Range("G36").Value = 16

$ws.Range("F37").Value = 48
$ws.Range("G37").Value = 16

$ws.Range("F38").Value = 8
$ws.Range("G38").Value = 16

# --- HUD: Grid / Health Bar Frame / Heath Bar / Reticule / Targeting Arrow ---
$ws.Range("F40").Value = 150
$ws.Range("G40").Value = 150
$ws.Range("I40").Value = "hud/grid"

$ws.Range("F41").Value = 36
$ws.Range("G41").Value = 5

$ws.Range("F42").Value = 26
$ws.Range("G42").Value = 3

$ws.Range("F43").Value = 24
$ws.Range("G43").Value = 24

$ws.Range("F44").Value = 12
$ws.Range("G44").Value = 6
$ws.Range("I44").Value = "hud/targeting-arrow"

# --- New "Start Screen" sprites ---
$ws.Range("B46").Value = "Logo"
$ws.Range("C46").Value = "Start Screen"
$ws.Range("D46").Value = 1
$ws.Range("E46").Value = "no"
$ws.Range("F46").Value = 390
$ws.Range("G46").Value = 164
$ws.Range("H46").Value = "hud_v1"
$ws.Range("I46").Value = "hud/logo"

$ws.Range("B47").Value = "Play Over"
$ws.Range("C47").Value = "Start Screen"
$ws.Range("D47").Value = 1
$ws.Range("E47").Value = "no"
$ws.Range("F47").Value = 168
$ws.Range("G47").Value = 72
$ws.Range("H47").Value = "hud_v1"
$ws.Range("I47").Value = "hud/play-over"

$ws.Range("B48").Value = "Play Out"
$ws.Range("C48").Value = "Start Screen"
$ws.Range("D48").Value = 1
$ws.Range("E48").Value = "no"
$ws.Range("F48").Value = 168
$ws.Range("G48").Value = 72
$ws.Range("H48").Value = "hud_v1"
$ws.Range("I48").Value = "hud/play-out"

$ws.Range("B49").Value = "Play Down"
$ws.Range("C49").Value = "Start Screen"
$ws.Range("D49").Value = 1
$ws.Range("E49").Value = "no"
$ws.Range("F49").Value = 168
$ws.Range("G49").Value = 72
$ws.Range("H49").Value = "hud_v1"
$ws.Range("I49").Value = "hud/play-down"

$ws.Range("B50").Value = "Play Up"
$ws.Range("C50").Value = "Start Screen"
$ws.Range("D50").Value = 1
$ws.Range("E50").Value = "no"
$ws.Range("F50").Value = 168
$ws.Range("G50").Value = 72
$ws.Range("H50").Value = "hud_v1"
$ws.Range("I50").Value = "hud/play-up"

# --- Update the selection / scrolled view to match the new active cell ---
$ws.Range("G42").Select()
